$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("detailed")

# Rows 4-9 in column F ("Difference") currently hold numeric 0.0 values
# comparing Source1/Source2 for NUM_INT, NUM_FLOAT, NUM_DOUBLE, NUM_DECIMAL1-3.
# They should instead report the string "no" (same as used elsewhere in the
# sheet), matching the fix for null-valued join columns.
foreach ($row in 4..9) {
    $ws.Cells.Item($row, 6).Value = "no"
}
